# Update TPM-derived values in the worksheet per new TPM calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.789499
$ws.Range("H2").Value = 5.368497
$ws.Range("I2").Value = 0.01244533957901722
$ws.Range("J2").Value = 0.01244533957901722
$ws.Range("M2").Value = 2.209083333333334
$ws.Range("N2").Value = 6.62725
$ws.Range("O2").Value = 0.2666231470852957
$ws.Range("P2").Value = 0.2666231470852958
$ws.Range("Q2").Value = 3.953152415916667
$ws.Range("R2").Value = 35.57837174325
$ws.Range("S2").Value = 0.00331821560510276
$ws.Range("T2").Value = 0.003318215605102761
$ws.Range("G3").Value = 1.789499
$ws.Range("H3").Value = 5.368497
$ws.Range("I3").Value = 0.01244533957901722
$ws.Range("J3").Value = 0.01244533957901722
$ws.Range("O3").Value = 0.07376078220023909
$ws.Range("P3").Value = 0.0737607822002391
$ws.Range("Q3").Value = 1.093632032861333
$ws.Range("R3").Value = 9.842688295752
$ws.Range("S3").Value = 0.0009179779820959043
$ws.Range("T3").Value = 0.0009179779820959044
$ws.Range("G4").Value = 1.789499
$ws.Range("H4").Value = 5.368497
$ws.Range("I4").Value = 0.01244533957901722
$ws.Range("J4").Value = 0.01244533957901722
$ws.Range("O4").Value = 0.659616070714465
$ws.Range("P4").Value = 0.6596160707144652
$ws.Range("Q4").Value = 9.779956811807333
$ws.Range("R4").Value = 88.019611306266
$ws.Range("S4").Value = 0.008209145991818552
$ws.Range("T4").Value = 0.008209145991818554
$ws.Range("H5").Value = 311.722962
$ws.Range("I5").Value = 0.7226413867171911
$ws.Range("J5").Value = 0.7226413867171912
$ws.Range("M5").Value = 2.209083333333334
$ws.Range("N5").Value = 6.62725
$ws.Range("O5").Value = 0.2666231470852957
$ws.Range("P5").Value = 0.2666231470852958
$ws.Range("Q5").Value = 229.5406666571667
$ws.Range("R5").Value = 2065.8659999145
$ws.Range("S5").Value = 0.1926729207406197
$ws.Range("T5").Value = 0.1926729207406198
$ws.Range("H6").Value = 311.722962
$ws.Range("I6").Value = 0.7226413867171911
$ws.Range("J6").Value = 0.7226413867171912
$ws.Range("O6").Value = 0.07376078220023909
$ws.Range("P6").Value = 0.0737607822002391
$ws.Range("S6").Value = 0.05330259393452548
$ws.Range("T6").Value = 0.0533025939345255
$ws.Range("H7").Value = 311.722962
$ws.Range("I7").Value = 0.7226413867171911
$ws.Range("J7").Value = 0.7226413867171912
$ws.Range("O7").Value = 0.659616070714465
$ws.Range("P7").Value = 0.6596160707144652
$ws.Range("R7").Value = 5110.878137862036
$ws.Range("S7").Value = 0.4766658720420458
$ws.Range("T7").Value = 0.476665872042046
$ws.Range("I8").Value = 0.2649132737037916
$ws.Range("J8").Value = 0.2649132737037916
$ws.Range("M8").Value = 2.209083333333334
$ws.Range("N8").Value = 6.62725
$ws.Range("O8").Value = 0.2666231470852957
$ws.Range("P8").Value = 0.2666231470852958
$ws.Range("Q8").Value = 84.14736627325
$ws.Range("R8").Value = 757.32629645925
$ws.Range("S8").Value = 0.07063201073957323
$ws.Range("T8").Value = 0.07063201073957324
$ws.Range("I9").Value = 0.2649132737037916
$ws.Range("J9").Value = 0.2649132737037916
$ws.Range("O9").Value = 0.07376078220023909
$ws.Range("P9").Value = 0.0737607822002391
$ws.Range("Q9").Value = 23.27920746663199
$ws.Range("S9").Value = 0.01954021028361769
$ws.Range("T9").Value = 0.0195402102836177
$ws.Range("I10").Value = 0.2649132737037916
$ws.Range("J10").Value = 0.2649132737037916
$ws.Range("O10").Value = 0.659616070714465
$ws.Range("P10").Value = 0.6596160707144652
$ws.Range("S10").Value = 0.1747410526806006
$ws.Range("T10").Value = 0.1747410526806006
